$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 19; B = "passed" },
    @{ Row = 20; B = "failed" },
    @{ Row = 21; B = "failed" },
    @{ Row = 22; B = "failed" },
    @{ Row = 23; B = "failed" },
    @{ Row = 24; B = "passed" },
    @{ Row = 25; B = "passed" },
    @{ Row = 26; B = "passed" },
    @{ Row = 27; B = "failed" }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row
    $ws.Cells.Item($rowIndex, 1).Value = "test-with-jdbc;country-testing-with-database"
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Borders.LineStyle = -4142
    $ws.Cells.Item($rowIndex, 4).Value = "28.07.21"
}
